$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.226.90"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.443.71"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9103"
$ws.Range("E5").Value = "  -9.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.83"
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3660"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3126"
$ws.Range("E8").Value = "  +3.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.11"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.021"
$ws.Range("E10").Value = "  +5.89%  "
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.389"
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.66"
$ws.Range("E14").Value = "  +6.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.072"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").Value = "1.441.25"
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9392"
$ws.Range("E18").Value = "  -6.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05629"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.67"
$ws.Range("E20").Value = "  -5.01%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.417"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.44"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.82"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.246"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "20.214.70"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.179"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "137.57"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.96"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").Value = "1.594.33"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.12"
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.776"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8039"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.839"
$ws.Range("E33").Value = "  -7.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07698"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06004"
$ws.Range("E35").Value = "  +6.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.451"
$ws.Range("E36").Value = "  +9.25%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.687"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.142"
$ws.Range("E38").Value = "  +10.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01990"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.17"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9291"
$ws.Range("E41").Value = "  -7.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1835"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.056"
$ws.Range("E43").Value = "  -14.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.521"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5237"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.27"
$ws.Range("E47").Value = "  +9.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5145"
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.766"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06333"
$ws.Range("E50").Value = "  +3.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9907"
$ws.Range("E51").Value = "  -1.14%  "
